$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 215
$ws.Range("B1").Value = 261.3999999999996
$ws.Range("C1").Value = 328.3999999999996

$ws.Range("A2").Value = 215
$ws.Range("B2").Value = 292.7999999999994
$ws.Range("C2").Value = 617.5499999565909
